$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values scraped from the updated "pl_mw.xlsx" case results (380 kV case)
$updates = @{
    "B2" = 0.7400215443108777
    "C2" = 0.06991220709525692
    "D2" = 0.02506600627322797
    "F2" = 5.700814329882093
    "G2" = 0.002622067671864311
    "J2" = 0.2942163618212561
    "K2" = 0.7310898725474715
    "M2" = 0.313472079305825
    "B3" = 0.7226644983766448
    "C3" = 0.06891224525639927
    "D3" = 0.02470331618800969
    "F3" = 5.50012605520547
    "G3" = 0.002627042174946926
    "J3" = 0.2895808351353821
    "K3" = 0.7151414572838917
    "M3" = 0.3117306321930222
    "B4" = 0.7126990819151899
    "C4" = 0.06840214110688692
    "D4" = 0.024535234021279
    "F4" = 5.377380857029863
    "G4" = 0.002630254225390262
    "J4" = 0.2868389137293406
    "K4" = 0.7060995805327224
    "M4" = 0.3109649248613024
    "B5" = 0.7088116744787953
    "C5" = 0.06822026226190303
    "D5" = 0.02448031369384296
    "F5" = 5.327477556254024
    "G5" = 0.002631602952586509
    "J5" = 0.2857477154179335
    "K5" = 0.7026031470305014
    "M5" = 0.3107291120810345
    "B6" = 0.7081766458878747
    "C6" = 0.06819162761671294
    "D6" = 0.02447200947561612
    "F6" = 5.319198066501741
    "G6" = 0.00263182931499454
    "J6" = 0.2855681004982813
    "K6" = 0.7020339191212059
    "M6" = 0.3106945558536616
    "B7" = 0.7126459526814415
    "C7" = 0.06839958314954231
    "D7" = 0.02453443859241133
    "F7" = 5.376707377655322
    "G7" = 0.002630272253358513
    "J7" = 0.2868240916315941
    "K7" = 0.7060516650570889
    "M7" = 0.3109614361449431
    "B8" = 0.7338929837803221
    "C8" = 0.06954578244769039
    "D8" = 0.02492952082802447
    "F8" = 5.631515351123824
    "G8" = 0.002623750237117082
    "J8" = 0.292596336267863
    "K8" = 0.7254347657665221
    "M8" = 0.3128085571751313
    "B9" = 0.7810702839823591
    "C9" = 0.07262386028706658
    "D9" = 0.02614481582664752
    "F9" = 6.135188310761436
    "G9" = 0.002612205363109728
    "J9" = 0.304747523426613
    "K9" = 0.7694275219746203
    "M9" = 0.3188453753830487
    "B10" = 0.8191284947822908
    "C10" = 0.07540068205628359
    "D10" = 0.02731641654120409
    "F10" = 6.507987489308078
    "G10" = 0.002604473235312102
    "J10" = 0.3141891724333732
    "K10" = 0.8054403338308873
    "M10" = 0.324762711178181
    "B11" = 0.8371883230258561
    "C11" = 0.07677788037881328
    "D11" = 0.0279121429625917
    "F11" = 6.67825145877282
    "G11" = 0.002601116606448463
    "J11" = 0.3185977053946516
    "K11" = 0.8226351735608262
    "M11" = 0.3277788342285604
    "B12" = 0.8441351235752848
    "C12" = 0.07731594506769568
    "D12" = 0.02814693283512071
    "F12" = 6.742828176019657
    "G12" = 0.002599868508248739
    "J12" = 0.3202835315121888
    "K12" = 0.8292639784884557
    "M12" = 0.3289677684637056
    "B13" = 0.8426341968151121
    "C13" = 0.07719932465680301
    "D13" = 0.0280959545832502
    "F13" = 6.728915843959953
    "G13" = 0.002600136288218366
    "J13" = 0.319919727346516
    "K13" = 0.8278311104571685
    "M13" = 0.3287096264209808
    "B14" = 0.8377576742444717
    "C14" = 0.0768218147059514
    "D14" = 0.02793127384460092
    "F14" = 6.683562166325999
    "G14" = 0.002601013464587443
    "J14" = 0.3187360698991881
    "K14" = 0.8231781707372647
    "M14" = 0.3278757098066833
    "B15" = 0.8347847355639715
    "C15" = 0.07659273886422113
    "D15" = 0.02783160566811205
    "F15" = 6.655795068470752
    "G15" = 0.002601553749768073
    "J15" = 0.3180131854268495
    "K15" = 0.8203434318220388
    "M15" = 0.3273710105987675
    "B16" = 0.8179633176756624
    "C16" = 0.07531298625363547
    "D16" = 0.02727876235142901
    "F16" = 6.496874300494852
    "G16" = 0.002604695822303427
    "J16" = 0.3139033557548032
    "K16" = 0.8043330153238912
    "M16" = 0.324572139717425
    "B17" = 0.8078356044377699
    "C17" = 0.07455720929154097
    "D17" = 0.02695581018733861
    "F17" = 6.399557657238518
    "G17" = 0.002606664461130236
    "J17" = 0.311411237976003
    "K17" = 0.7947196947414739
    "M17" = 0.3229382957831319
    "B18" = 0.8020806865741292
    "C18" = 0.07413322789577137
    "D18" = 0.02677595321775783
    "F18" = 6.343647044056553
    "G18" = 0.002607811909225894
    "J18" = 0.3099885133083689
    "K18" = 0.7892667960275901
    "M18" = 0.3220290664535455
    "B19" = 0.8001442243244412
    "C19" = 0.07399151182653441
    "D19" = 0.02671606412651784
    "F19" = 6.324727421525154
    "G19" = 0.002608203019443249
    "J19" = 0.309508633953314
    "K19" = 0.7874336449601174
    "M19" = 0.3217264531173143
    "B20" = 0.8089064398795358
    "C20" = 0.0746365522905279
    "D20" = 0.02698957744942021
    "F20" = 6.40991059320109
    "G20" = 0.002606453330425839
    "J20" = 0.3116754226241909
    "K20" = 0.7957351342626282
    "M20" = 0.323109061988049
    "B21" = 0.8391870935419661
    "C21" = 0.07693224812605592
    "D21" = 0.02797939343769684
    "F21" = 6.696880853237303
    "G21" = 0.002600755194130097
    "J21" = 0.3190832926644873
    "K21" = 0.8245416578244829
    "M21" = 0.3281193800530673
    "B22" = 0.8596066360339591
    "C22" = 0.07852915587297105
    "D22" = 0.02868001002023846
    "F22" = 6.885025843900564
    "G22" = 0.002597165042162771
    "J22" = 0.3240204606677679
    "K22" = 0.8440535424253994
    "M22" = 0.3316667414640904
    "B23" = 0.8486505760978389
    "C23" = 0.07766797083169763
    "D23" = 0.02830110418145182
    "F23" = 6.784553603937468
    "G23" = 0.002599068964474509
    "J23" = 0.3213766138366765
    "K23" = 0.8335767661375826
    "M23" = 0.329748429657883
    "B24" = 0.8084221045273523
    "C24" = 0.07460064854845427
    "D24" = 0.02697429317877464
    "F24" = 6.405229912734541
    "G24" = 0.002606548733927934
    "J24" = 0.3115559535215056
    "K24" = 0.7952758238305648
    "M24" = 0.3230317648741163
    "B25" = 0.7677137287984976
    "C25" = 0.07170132632137438
    "D25" = 0.02576783797189108
    "F25" = 5.998471334376376
    "G25" = 0.002615196223590588
    "J25" = 0.3013705302039966
    "K25" = 0.7568813310563201
    "M25" = 0.3189527323609742
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
